$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01555952548175
$ws.Range("C2").Value = 0.01555952548175
$ws.Range("D2").Value = 0.0111352794092782
$ws.Range("E2").Value = 0.000076894587676767
$ws.Range("F2").Value = 0.9098

# Row 3
$ws.Range("B3").Value = 0.564522505049724
$ws.Range("C3").Value = 0.564522505049724
$ws.Range("D3").Value = 0.404004340230516
$ws.Range("E3").Value = 0.00278984891351403
$ws.Range("F3").Value = 0.5262

# Row 4
$ws.Range("B4").Value = 0.554910439568874
$ws.Range("C4").Value = 0.554910439568874
$ws.Range("D4").Value = 0.397125400705331
$ws.Range("E4").Value = 0.00274234644869022
$ws.Range("F4").Value = 0.5316

# Row 5
$ws.Range("B5").Value = 201.21378072517
$ws.Range("C5").Value = 1.39731792170257
$ws.Range("E5").Value = 0.994390910050119

# Row 6
$ws.Range("B6").Value = 202.348773195271
